$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(22, 1).Value = -3.719183167544286
$ws.Cells.Item(22, 2).Value = 7.010346993533103
$ws.Cells.Item(22, 3).Value = 0.6059804959730433

$ws.Cells.Item(23, 1).Value = -3.724935540285977
$ws.Cells.Item(23, 2).Value = 10.48463944521817
$ws.Cells.Item(23, 3).Value = -5.378288039294167

$ws.Cells.Item(24, 1).Value = 2.686929789456455
$ws.Cells.Item(24, 2).Value = 0.002814553000709097
$ws.Cells.Item(24, 3).Value = -3.492711760781028

$ws.Cells.Item(25, 1).Value = 3.455663386258163
$ws.Cells.Item(25, 2).Value = -6.499211259321735
$ws.Cells.Item(25, 3).Value = -1.833242598446949

$ws.Cells.Item(26, 1).Value = 0.03122558593747993
$ws.Cells.Item(26, 2).Value = -2.420187681913326
$ws.Cells.Item(26, 3).Value = 0.5764146804809926

$ws.Cells.Item(27, 1).Value = -7.196380597894843
$ws.Cells.Item(27, 2).Value = -4.824799558791248
$ws.Cells.Item(27, 3).Value = 5.860242297432641

$ws.Cells.Item(28, 1).Value = -1.569692446968803
$ws.Cells.Item(28, 2).Value = -11.00871762362393
$ws.Cells.Item(28, 3).Value = 5.813380349766132

$ws.Cells.Item(29, 1).Value = 3.11778094551775
$ws.Cells.Item(29, 2).Value = -8.70344656163989
$ws.Cells.Item(29, 3).Value = 3.647271784869103

$ws.Cells.Item(30, 1).Value = 4.020912179079962
$ws.Cells.Item(30, 2).Value = 1.777516035600222
$ws.Cells.Item(30, 3).Value = -0.3423178889535121

$ws.Cells.Item(31, 1).Value = 2.335198922590785
$ws.Cells.Item(31, 2).Value = 4.166226341507659
$ws.Cells.Item(31, 3).Value = 0.4891112804412763
